$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 51.12915033333334
$ws.Range("H2").Value = 153.387451
$ws.Range("I2").Value = 0.8013178159252168
$ws.Range("J2").Value = 0.8013178159252169
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 35.48871133333333
$ws.Range("N2").Value = 106.466134
$ws.Range("O2").Value = 0.7670904531193554
$ws.Range("P2").Value = 0.7670904531193554
$ws.Range("Q2").Value = 1814.507656898271
$ws.Range("R2").Value = 16330.56891208444
$ws.Range("S2").Value = 0.6146832465106867
$ws.Range("T2").Value = 0.6146832465106868

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 51.12915033333334
$ws.Range("H3").Value = 153.387451
$ws.Range("I3").Value = 0.8013178159252168
$ws.Range("J3").Value = 0.8013178159252169
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.613261333333333
$ws.Range("N3").Value = 7.839784
$ws.Range("O3").Value = 0.05648578787427251
$ws.Range("P3").Value = 0.0564857878742725
$ws.Range("Q3").Value = 133.6138315722871
$ws.Range("R3").Value = 1202.524484150584
$ws.Range("S3").Value = 0.04526306817022714
$ws.Range("T3").Value = 0.04526306817022714

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 51.12915033333334
$ws.Range("H4").Value = 153.387451
$ws.Range("I4").Value = 0.8013178159252168
$ws.Range("J4").Value = 0.8013178159252169
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.596082
$ws.Range("N4").Value = 1.788246
$ws.Range("O4").Value = 0.01288434531142903
$ws.Range("P4").Value = 0.01288434531142903
$ws.Range("Q4").Value = 30.477166188994
$ws.Range("R4").Value = 274.294495700946
$ws.Range("S4").Value = 0.01032445544458062
$ws.Range("T4").Value = 0.01032445544458062

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 51.12915033333334
$ws.Range("H5").Value = 153.387451
$ws.Range("I5").Value = 0.8013178159252168
$ws.Range("J5").Value = 0.8013178159252169
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.482776
$ws.Range("N5").Value = 1.448328
$ws.Range("O5").Value = 0.01043522987117622
$ws.Range("P5").Value = 0.01043522987117622
$ws.Range("Q5").Value = 24.68392668132534
$ws.Range("R5").Value = 222.155340131928
$ws.Range("S5").Value = 0.008361935609048507
$ws.Range("T5").Value = 0.008361935609048507

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 51.12915033333334
$ws.Range("H6").Value = 153.387451
$ws.Range("I6").Value = 0.8013178159252168
$ws.Range("J6").Value = 0.8013178159252169
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.083219666666667
$ws.Range("N6").Value = 21.249659
$ws.Range("O6").Value = 0.1531041838237668
$ws.Range("P6").Value = 0.1531041838237668
$ws.Range("Q6").Value = 362.1590031810233
$ws.Range("R6").Value = 3259.431028629209
$ws.Range("S6").Value = 0.1226851101906737
$ws.Range("T6").Value = 0.1226851101906737

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.875170333333333
$ws.Range("H7").Value = 8.625510999999999
$ws.Range("I7").Value = 0.04506089377389114
$ws.Range("J7").Value = 0.04506089377389114
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 35.48871133333333
$ws.Range("N7").Value = 106.466134
$ws.Range("O7").Value = 0.7670904531193554
$ws.Range("P7").Value = 0.7670904531193554
$ws.Range("Q7").Value = 102.0360899938304
$ws.Range("R7").Value = 918.3248099444741
$ws.Range("S7").Value = 0.0345657814229773
$ws.Range("T7").Value = 0.0345657814229773

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.875170333333333
$ws.Range("H8").Value = 8.625510999999999
$ws.Range("I8").Value = 0.04506089377389114
$ws.Range("J8").Value = 0.04506089377389114
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.613261333333333
$ws.Range("N8").Value = 7.839784
$ws.Range("O8").Value = 0.05648578787427251
$ws.Range("P8").Value = 0.0564857878742725
$ws.Range("Q8").Value = 7.513571458847111
$ws.Range("R8").Value = 67.62214312962399
$ws.Range("S8").Value = 0.002545300087137142
$ws.Range("T8").Value = 0.002545300087137142

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.875170333333333
$ws.Range("H9").Value = 8.625510999999999
$ws.Range("I9").Value = 0.04506089377389114
$ws.Range("J9").Value = 0.04506089377389114
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.596082
$ws.Range("N9").Value = 1.788246
$ws.Range("O9").Value = 0.01288434531142903
$ws.Range("P9").Value = 0.01288434531142903
$ws.Range("Q9").Value = 1.713837282634
$ws.Range("R9").Value = 15.424535543706
$ws.Range("S9").Value = 0.0005805801154244358
$ws.Range("T9").Value = 0.0005805801154244359

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.875170333333333
$ws.Range("H10").Value = 8.625510999999999
$ws.Range("I10").Value = 0.04506089377389114
$ws.Range("J10").Value = 0.04506089377389114
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.482776
$ws.Range("N10").Value = 1.448328
$ws.Range("O10").Value = 0.01043522987117622
$ws.Range("P10").Value = 0.01043522987117622
$ws.Range("Q10").Value = 1.388063232845333
$ws.Range("R10").Value = 12.492569095608
$ws.Range("S10").Value = 0.0004702207847312072
$ws.Range("T10").Value = 0.0004702207847312072

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.875170333333333
$ws.Range("H11").Value = 8.625510999999999
$ws.Range("I11").Value = 0.04506089377389114
$ws.Range("J11").Value = 0.04506089377389114
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 7.083219666666667
$ws.Range("N11").Value = 21.249659
$ws.Range("O11").Value = 0.1531041838237668
$ws.Range("P11").Value = 0.1531041838237668
$ws.Range("Q11").Value = 20.36546305008322
$ws.Range("R11").Value = 183.289167450749
$ws.Range("S11").Value = 0.006899011363621058
$ws.Range("T11").Value = 0.006899011363621058

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 9.802011
$ws.Range("H12").Value = 29.406033
$ws.Range("I12").Value = 0.153621290300892
$ws.Range("J12").Value = 0.153621290300892
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 35.48871133333333
$ws.Range("N12").Value = 106.466134
$ws.Range("O12").Value = 0.7670904531193554
$ws.Range("P12").Value = 0.7670904531193554
$ws.Range("Q12").Value = 347.860738865158
$ws.Range("R12").Value = 3130.746649786422
$ws.Range("S12").Value = 0.1178414251856913
$ws.Range("T12").Value = 0.1178414251856913

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 9.802011
$ws.Range("H13").Value = 29.406033
$ws.Range("I13").Value = 0.153621290300892
$ws.Range("J13").Value = 0.153621290300892
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.613261333333333
$ws.Range("N13").Value = 7.839784
$ws.Range("O13").Value = 0.05648578787427251
$ws.Range("P13").Value = 0.0564857878742725
$ws.Range("Q13").Value = 25.615216335208
$ws.Range("R13").Value = 230.536947016872
$ws.Range("S13").Value = 0.008677419616908225
$ws.Range("T13").Value = 0.008677419616908223

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.802011
$ws.Range("H14").Value = 29.406033
$ws.Range("I14").Value = 0.153621290300892
$ws.Range("J14").Value = 0.153621290300892
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.596082
$ws.Range("N14").Value = 1.788246
$ws.Range("O14").Value = 0.01288434531142903
$ws.Range("P14").Value = 0.01288434531142903
$ws.Range("Q14").Value = 5.842802320902
$ws.Range("R14").Value = 52.585220888118
$ws.Range("S14").Value = 0.001979309751423976
$ws.Range("T14").Value = 0.001979309751423976

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.802011
$ws.Range("H15").Value = 29.406033
$ws.Range("I15").Value = 0.153621290300892
$ws.Range("J15").Value = 0.153621290300892
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.482776
$ws.Range("N15").Value = 1.448328
$ws.Range("O15").Value = 0.01043522987117622
$ws.Range("P15").Value = 0.01043522987117622
$ws.Range("Q15").Value = 4.732175662536
$ws.Range("R15").Value = 42.589580962824
$ws.Range("S15").Value = 0.001603073477396502
$ws.Range("T15").Value = 0.001603073477396501

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.802011
$ws.Range("H16").Value = 29.406033
$ws.Range("I16").Value = 0.153621290300892
$ws.Range("J16").Value = 0.153621290300892
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 7.083219666666667
$ws.Range("N16").Value = 21.249659
$ws.Range("O16").Value = 0.1531041838237668
$ws.Range("P16").Value = 0.1531041838237668
$ws.Range("Q16").Value = 69.42979708808301
$ws.Range("R16").Value = 624.868173792747
$ws.Range("S16").Value = 0.02352006226947202
$ws.Range("T16").Value = 0.02352006226947202
